# Refresh the cryptos table with the latest scraped prices / volume-1h percentages
# (GitHub Actions cron update). A couple of rows (Stacks/FirstDigitalUSD and
# InjectiveProtocol/Arweave) also swapped rank order, so Coin/Link are rewritten too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.755.53"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "'3.172.49"

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'613.29"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").Value = "'146.49"
$ws.Range("E6").Value = "  -1.85%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'3.167.12"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("E10").Value = "  -0.13%  "

$ws.Range("D11").Value = "'5.48"
$ws.Range("E11").Value = "  -2.57%  "

$ws.Range("D12").Value = "'0.476"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "'0.0000260"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "'35.93"
$ws.Range("E14").Value = "  -2.52%  "

$ws.Range("D15").Value = "'3.694.88"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("E16").Value = "  +2.99%  "

$ws.Range("D17").Value = "'64.748.07"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").Value = "'3.172.42"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").Value = "'6.91"
$ws.Range("E19").Value = "  -1.41%  "

$ws.Range("D20").Value = "'481.17"
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").Value = "'14.71"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'0.723"
$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").Value = "'7.94"
$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D24").Value = "'13.80"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("D25").Value = "'84.37"
$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'8.74"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("D28").Value = "'2.82"
$ws.Range("E28").Value = "  -4.48%  "

$ws.Range("E29").Value = "  +2.16%  "

$ws.Range("E30").Value = "  -3.54%  "

$ws.Range("D31").Value = "'2.12"
$ws.Range("E31").Value = "  -5.78%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.72"
$ws.Range("E33").Value = "  -0.83%  "

$ws.Range("D34").Value = "'26.70"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").Value = "'1.13"
$ws.Range("E35").Value = "  +1.67%  "

$ws.Range("D36").Value = "'0.0₃0794"
$ws.Range("E36").Value = "  +5.90%  "

$ws.Range("D37").Value = "'6.03"
$ws.Range("E37").Value = "  -1.51%  "

$ws.Range("E38").Value = "  -2.16%  "

$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").Value = "'464.04"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("D41").Value = "'0.0401"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").Value = "'0.120"
$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("D43").Value = "'8.38"
$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").Value = "'2.862.27"
$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").Value = "'0.269"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = "  +4.74%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'26.79"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'36.59"
$ws.Range("E49").Value = "  +8.55%  "

$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("E51").Value = "  -0.45%  "
